$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")
$ws.Range("F2").Value = 3
$ws.Range("F3").Select()
